$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New user rows to append: username, password, role
$newRows = @(
    @("T1", "RCSS", "teacher"),
    @("T2", "RCSS", "teacher"),
    @("T3", "RCSS", "teacher"),
    @("T4", "RCSS", "teacher"),
    @("T5", "RCSS", "teacher")
)

$startRow = 3
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $newRows[$i][0]
    $ws.Cells.Item($r, 2).Value = $newRows[$i][1]
    $ws.Cells.Item($r, 3).Value = $newRows[$i][2]
}

# Resize the table (ListObject) to cover the new data range
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:C7"))

# Update selection to match the last-edited cell
$ws.Range("C7").Select()
